$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1048
$ws1.Range("F5").Value = 2856
$ws1.Range("F7").Value = 249
$ws1.Range("F8").Value = 23
$ws1.Range("F10").Value = 89
$ws1.Range("F11").Value = 114
$ws1.Range("F12").Value = 27
$ws1.Range("F13").Value = 2685
$ws1.Range("F14").Value = 910

# Sheet "全部类型" (All Types) - column F ("想去人数") updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1048
$ws4.Range("F6").Value = 2856
$ws4.Range("F8").Value = 249
$ws4.Range("F9").Value = 23
$ws4.Range("F12").Value = 89
$ws4.Range("F13").Value = 114
$ws4.Range("F14").Value = 27
$ws4.Range("F15").Value = 2685
$ws4.Range("F16").Value = 910
